$d = $word.ActiveDocument

# =====================================================================
# Edit 1: "... cloudy (r=.04)" -> "... cloudy (r=.02)"
#   The trailing "4)" becomes "2)", split into three runs (digit / close
#   paren) exactly as Word does when you retype a single character deep
#   inside an existing sentence, and the special "_GoBack" bookmark is
#   re-added at the edit point (which also removes it from its old spot).
# =====================================================================

$rng1 = $d.Content
$null = $rng1.Find.Execute("cloudy (r=.04)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraEnd1 = $rng1.End

# Temporarily give the trailing digit and the closing paren distinct,
# mismatched colors so Word is forced to keep them as separate runs
# instead of silently re-merging them with their neighbours while we
# edit their text.
$digitRng1 = $d.Range($paraEnd1 - 2, $paraEnd1 - 1)
$digitRng1.Font.Color = 128

$closeRng1 = $d.Range($paraEnd1 - 1, $paraEnd1)
$closeRng1.Font.Color = 255

# Replace the digit "4" -> "2".
$digitRng1b = $d.Range($paraEnd1 - 2, $paraEnd1 - 1)
$digitRng1b.Text = "2"

# Round-trip the closing paren through a placeholder so it drops any
# stale run identity (rsid) it inherited, matching a genuinely new run.
$closeRng1b = $d.Range($paraEnd1 - 1, $paraEnd1)
$closeRng1b.Text = "X"
$closeRng1c = $d.Range($paraEnd1 - 1, $paraEnd1)
$closeRng1c.Text = ")"

# Restore the original (identical) black color now that the runs exist
# as separate elements.
$digitRng1c = $d.Range($paraEnd1 - 2, $paraEnd1 - 1)
$digitRng1c.Font.Color = 0
$closeRng1d = $d.Range($paraEnd1 - 1, $paraEnd1)
$closeRng1d.Font.Color = 0

# Re-insert the special "_GoBack" bookmark between the digit run and the
# closing-paren run - this both places it here and removes it from its
# previous location further down in the document.
$bmRng1 = $d.Range($paraEnd1 - 1, $paraEnd1 - 1)
$d.Bookmarks.Add("_GoBack", $bmRng1)

# =====================================================================
# Edit 2: "Latitude (r=-.02)." -> "Latitude (r=-.00)."
#   The trailing "2" becomes "0", and the remainder of the sentence is
#   split off into its own run (still starting with the closing paren).
# =====================================================================

$rng2 = $d.Content
$null = $rng2.Find.Execute("Latitude (r=-.02)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraEnd2 = $rng2.End

$paraCount2 = $d.Range(0, $paraEnd2).Paragraphs.Count
$paraObj2 = $d.Paragraphs.Item($paraCount2)
$fullParaEnd2 = $paraObj2.Range.End - 1

# Distinct colors again so nothing merges back together while we edit.
$digitRng2 = $d.Range($paraEnd2 - 2, $paraEnd2 - 1)
$digitRng2.Font.Color = 128

$remRng2 = $d.Range($paraEnd2 - 1, $fullParaEnd2)
$remRng2.Font.Color = 255

# Replace the digit "2" -> "0".
$digitRng2b = $d.Range($paraEnd2 - 2, $paraEnd2 - 1)
$digitRng2b.Text = "0"

# Round-trip the remainder of the sentence (closing paren onward)
# through a placeholder so it drops its stale run identity too.
$remRng2b = $d.Range($paraEnd2 - 1, $fullParaEnd2)
$remainderText = $remRng2b.Text
$remRng2b.Text = "PH"
$remRng2c = $d.Range($paraEnd2 - 1, $paraEnd2 - 1 + 2)
$remRng2c.Text = $remainderText

# Restore the original black color for both new runs.
$digitRng2c = $d.Range($paraEnd2 - 2, $paraEnd2 - 1)
$digitRng2c.Font.Color = 0
$remEndNow2 = $paraEnd2 - 1 + $remainderText.Length
$remRng2d = $d.Range($paraEnd2 - 1, $remEndNow2)
$remRng2d.Font.Color = 0

Write-Output "done"
